$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Cell Markers (1)"

# Delete row 13 (was blank), shifting rows 14:24 up to 13:23
$ws.Rows("13:13").Delete()

# Select the row that now occupies position 13 (mirrors the final UI selection)
$ws.Rows("13:13").Select()
